$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the existing jx:if condition cell (A6) lastCell reference from B6 to A6
$ws.Range("A6").Value = 'jx:if(condition="person.age < 18", lastCell="A6")'

# 2. Add the new "Address" row (row 10): label + placeholder
$ws.Range("A10").Value = "Address:"
$ws.Range("B10").Value = '${person.address.addressLine}'

# 3. Add the new jx:if condition row for address (row 9), with the same highlight fill as A6
$ws.Range("A9").Value = 'jx:if(condition="person.addressExists ", lastCell="A9")'
$ws.Range("A9").Interior.ColorIndex = $ws.Range("A6").Interior.ColorIndex

# 4. Duplicate the explanatory comment from A6 onto A9
$commentText = $ws.Range("A6").Comment.Text()
$ws.Range("A9").AddComment($commentText)

# 5. Widen column A to fit the new, longer content
$ws.Columns.Item(1).ColumnWidth = 43.43

# 6. Update the active selection to the new last cell
[void]$ws.Range("A10").Select()

Write-Output "done"
